# Updates the summary table to reflect the corrected TableComparator output
# (fix to the comparison logic + a new generic BigDecimal numeric type means
# numeric examples are now reported as bare ids, without their country
# suffixes, and several counts/percentages have changed).
#
# All the values in this table are stored as *text* in the workbook (not
# numbers), so each target range has its number format forced to "@" (text)
# before the value is written. This prevents Excel from reinterpreting
# strings like "10" or "14.3%" as numeric/percentage values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
}

# Row 2: REF unique id count 11 -> 10
Set-TextValue "D2" "10"

# Row 7: denominador for Quality global 11 -> 10 (follows REF unique count)
Set-TextValue "E7" "10"

# Row 8: MATCH / 1:1 (exact matches)
Set-TextValue "D8" "1"
Set-TextValue "E8" "7"
Set-TextValue "F8" "14.3%"
Set-TextValue "G8" "NULL"

# Row 9: NO MATCH / 1:1 (match not identical)
Set-TextValue "D9" "6"
Set-TextValue "E9" "7"
Set-TextValue "F9" "85.7%"
Set-TextValue "G9" "1,2,4,7,8,9"

# Row 10: GAP / 1:0 (only in reference)
Set-TextValue "D10" "3"
Set-TextValue "E10" "10"
Set-TextValue "F10" "30.0%"
Set-TextValue "G10" "10,3,5"

# Row 11: GAP / 0:1 (only in new)
Set-TextValue "D11" "1"
Set-TextValue "F11" "12.5%"
Set-TextValue "G11" "6"

# Row 12: DUPS / duplicates (both)
Set-TextValue "D12" "2"
Set-TextValue "E12" "7"
Set-TextValue "F12" "28.6%"
Set-TextValue "G12" "4,NULL"

# Row 13: DUPS / duplicates (ref)
Set-TextValue "E13" "10"
Set-TextValue "F13" "10.0%"
Set-TextValue "G13" "5"

# Row 14: DUPS / duplicates (new)
Set-TextValue "D14" "1"
Set-TextValue "F14" "12.5%"
Set-TextValue "G14" "6"

$wb.Save()
